$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80, shifting existing rows 80-120 down to 81-121.
$ws.Rows("80:80").Insert()

# Populate the newly inserted row 80 with this week's new data point
# (same market/product/quality as the row above, new date & prices).
$ws.Range("A80").Value = 4
$ws.Range("B80").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C80").Value = "Los Lagos"
$ws.Range("D80").Value = 44992
$ws.Range("E80").Value = 10
$ws.Range("F80").Value = "Fruta"
$ws.Range("G80").Value = 100104
$ws.Range("H80").Value = "Frutos de pepita"
$ws.Range("I80").Value = 100104003
$ws.Range("J80").Value = "Membrillo"
$ws.Range("K80").Value = "Champion"
$ws.Range("L80").Value = "Primera"
$ws.Range("M80").Value = 300
$ws.Range("N80").Value = 17000
$ws.Range("O80").Value = 18000
$ws.Range("P80").Value = 17500
$ws.Range("Q80").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R80").Value = "Región de O'Higgins"
$ws.Range("S80").Value = 972
$ws.Range("T80").Value = 18
